# Applies the "fixed html and php slides" commit to the presentation.
#
# Slide 12 ("OUR SCHEDULE"): update two schedule line items' hours/sessions.
# Slide 2 ("BEFORE WE GET STARTED"): resize the checklist textbox, tweak
# several checklist bullet lines, and trim the browser-extensions bullet.
#
# NOTE: this runtime recomputes a shape's Height as a side effect of any
# TextRange.Text write on a shape that uses <a:spAutoFit/> (even a no-op
# write). Both edited shapes below use spAutoFit, so every text edit is
# done first and the shape's position/size is (re)asserted afterwards to
# land on the exact EMU values from the target deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 : "OUR SCHEDULE" -- TextBox 4 (shape 3)
# ---------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$schedule = $slide12.Shapes.Item(3)
$scheduleText = $schedule.TextFrame.TextRange

# " HTML & CSS & Bootstrap 4 : 16 hours - 4 sessions." -> 28 hours - 7 sessions.
$htmlPara = $scheduleText.Paragraphs(3, 1)
$htmlRun = $htmlPara.Runs(1, 1)
$htmlRun.Text = " HTML & CSS & Bootstrap 4 : 28 hours " + [char]0x2013 + " 7 sessions."

# "Project News Website Template : 32 hours - 7 sessions." -> 16 hours - 4 sessions.
$projPara = $scheduleText.Paragraphs(5, 1)
$projRun = $projPara.Runs(1, 1)
$projRun.Text = "Project News Website Template : 16 hours " + [char]0x2013 + " 4 sessions."

# The edits above re-trigger spAutoFit on this shape even though the
# rendered line count is unchanged; pin the height back to its original
# (unchanged) size so the shape's xfrm stays untouched, as in the target.
$schedule.Height = 288.7021484375

# ---------------------------------------------------------------------
# Slide 2 : "BEFORE WE GET STARTED" -- TextBox 3 (shape 3)
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$checklist = $slide2.Shapes.Item(3)
$checklistText = $checklist.TextFrame.TextRange

# Paragraph 1 : "Using Windows." -> "Using Windows 10 " + "or Ubuntu 16.04."
$winPara = $checklistText.Paragraphs(1, 1)
$winRun = $winPara.Runs(0, 1)
$newWinRun = $winRun.InsertBefore("Using Windows 10 ")
$winRun.Text = "or Ubuntu 16.04."

# Paragraph 2 : "Install Git." -> "Installed Git."
$gitPara = $checklistText.Paragraphs(2, 1)
$gitRun = $gitPara.Runs(1, 1)
$gitRun.Text = "Installed Git."

# Paragraph 3 : VS Code extensions list rewrite.
$vscodePara = $checklistText.Paragraphs(3, 1)
$vscodeRun = $vscodePara.Runs(1, 1)
$vscodeRun.Text = "Installed Visual Studio Code " + [char]0x2013 + " with extensions with extensions : Live server, Prettier, HTML CSS Support, IntelliSense for CSS classes, Auto Close Tag, HTML Snippets, Bootstrap 4 Snippets, Material Theme, Material Icons " + [char]0x2026

# Paragraph 4 : collapse " Translate, Color " + "Zilla" runs into " Translate."
$extPara = $checklistText.Paragraphs(4, 1)
$translateRun = $extPara.Runs(3, 1)
$translateRun.Text = " Translate."
$zillaRun = $extPara.Runs(4, 1)
$zillaRun.Text = ""

# Grow the textbox to fit the extra line of text (set last, since editing
# the runs above recomputes -- and would otherwise clobber -- the height).
$checklist.Top = 126.10929870605469
$checklist.Left = 52.4669303894043
$checklist.Width = 855.0660400390625
$checklist.Height = 295.20001220703125
